$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-01-21 02:04:16"

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
